# Apply the betting-odds update for 2025-11-25 workbook.
# Updates Time for row 2 and numeric odds/prices across rows 2-23.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "00:10:00"
$ws.Range("F2").Value = 1.69
$ws.Range("G2").Value = 1.71
$ws.Range("I2").Value = 5.7
$ws.Range("J2").Value = 4.3
$ws.Range("K2").Value = 4.5
$ws.Range("L2").Value = 1.29
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 5.1
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 2.42
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 1.57
$ws.Range("S2").Value = 2.66
$ws.Range("T2").Value = 1.71
$ws.Range("U2").Value = 2.34
$ws.Range("W2").Value = 2.32
$ws.Range("X2").Value = 21
$ws.Range("Y2").Value = 29
$ws.Range("Z2").Value = 55
$ws.Range("AA2").Value = 150
$ws.Range("AC2").Value = 10
$ws.Range("AD2").Value = 22
$ws.Range("AE2").Value = 80
$ws.Range("AF2").Value = 12
$ws.Range("AG2").Value = 10.5
$ws.Range("AH2").Value = 19
$ws.Range("AI2").Value = 75
$ws.Range("AJ2").Value = 18
$ws.Range("AK2").Value = 16.5
$ws.Range("AL2").Value = 36
$ws.Range("AM2").Value = 110
$ws.Range("AN2").Value = 8.199999999999999
$ws.Range("AO2").Value = 70
# Row 3
$ws.Range("Q3").Value = 2
$ws.Range("S3").Value = 3.65
$ws.Range("T3").Value = 2.34
$ws.Range("U3").Value = 1.64
$ws.Range("W3").Value = 3.4
$ws.Range("Y3").Value = 80
$ws.Range("AB3").Value = 6.6
$ws.Range("AK3").Value = 18.5
$ws.Range("AL3").Value = 150
$ws.Range("AN3").Value = 8
# Row 4
$ws.Range("H4").Value = 1.96
$ws.Range("N4").Value = 4.8
$ws.Range("P4").Value = 2.26
$ws.Range("Q4").Value = 1.76
$ws.Range("R4").Value = 1.5
$ws.Range("S4").Value = 2.9
$ws.Range("AE4").Value = 17.5
$ws.Range("AG4").Value = 16
$ws.Range("AI4").Value = 29
$ws.Range("AO4").Value = 10.5
# Row 5
$ws.Range("F5").Value = 1.74
$ws.Range("AC5").Value = 10.5
$ws.Range("AH5").Value = 15.5
# Row 6
$ws.Range("G6").Value = 1.47
$ws.Range("J6").Value = 5.5
$ws.Range("U6").Value = 2.68
# Row 7
$ws.Range("H7").Value = 3.35
$ws.Range("P7").Value = 2.86
# Row 8
$ws.Range("F8").Value = 1.92
$ws.Range("L8").Value = 1.24
$ws.Range("AJ8").Value = 65
# Row 9
$ws.Range("S9").Value = 2.04
$ws.Range("W9").Value = 3.15
$ws.Range("X9").Value = 42
$ws.Range("AA9").Value = 230
$ws.Range("AB9").Value = 14.5
$ws.Range("AC9").Value = 14.5
$ws.Range("AF9").Value = 12
$ws.Range("AG9").Value = 11.5
$ws.Range("AJ9").Value = 14
$ws.Range("AK9").Value = 14
$ws.Range("AN9").Value = 4.4
# Row 10
$ws.Range("F10").Value = 2.2
$ws.Range("G10").Value = 2.34
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 4.5
$ws.Range("L10").Value = 1.57
$ws.Range("M10").Value = 1.12
$ws.Range("N10").Value = 2.68
$ws.Range("P10").Value = 1.55
$ws.Range("R10").Value = 1.2
$ws.Range("S10").Value = 5.3
$ws.Range("W10").Value = 1.74
$ws.Range("X10").Value = 8.800000000000001
$ws.Range("Y10").Value = 11.5
$ws.Range("Z10").Value = 120
$ws.Range("AD10").Value = 18.5
$ws.Range("AF10").Value = 13
$ws.Range("AH10").Value = 65
$ws.Range("AJ10").Value = 150
$ws.Range("AK10").Value = 80
$ws.Range("AL10").Value = 480
$ws.Range("AM10").Value = 220
$ws.Range("AO10").Value = 120
# Row 11
$ws.Range("F11").Value = 2.84
$ws.Range("P11").Value = 2.1
$ws.Range("U11").Value = 2.38
$ws.Range("AB11").Value = 28
$ws.Range("AE11").Value = 80
$ws.Range("AF11").Value = 22
$ws.Range("AI11").Value = 150
# Row 12
$ws.Range("F12").Value = 4
$ws.Range("J12").Value = 3.4
$ws.Range("K12").Value = 3.95
$ws.Range("L12").Value = 1.35
$ws.Range("X12").Value = 980
$ws.Range("Z12").Value = 980
$ws.Range("AA12").Value = 980
$ws.Range("AB12").Value = 980
$ws.Range("AF12").Value = 980
$ws.Range("AG12").Value = 980
$ws.Range("AH12").Value = 980
$ws.Range("AI12").Value = 980
# Row 13
$ws.Range("F13").Value = 1.88
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 4.6
$ws.Range("I13").Value = 5.6
$ws.Range("N13").Value = 2.92
$ws.Range("V13").Value = 1.21
$ws.Range("W13").Value = 2
$ws.Range("Y13").Value = 1000
$ws.Range("AB13").Value = 14
$ws.Range("AC13").Value = 15
$ws.Range("AF13").Value = 1000
$ws.Range("AL13").Value = 130
# Row 14
$ws.Range("G14").Value = 1.95
$ws.Range("N14").Value = 5
$ws.Range("O14").Value = 1.23
$ws.Range("W14").Value = 2.04
$ws.Range("Y14").Value = 20
$ws.Range("AC14").Value = 9
$ws.Range("AG14").Value = 10.5
$ws.Range("AJ14").Value = 22
$ws.Range("AK14").Value = 17
# Row 15
$ws.Range("F15").Value = 3.35
$ws.Range("G15").Value = 3.4
$ws.Range("V15").Value = 1.8
$ws.Range("AL15").Value = 36
# Row 16
$ws.Range("F16").Value = 2.36
$ws.Range("G16").Value = 2.38
$ws.Range("H16").Value = 2.96
$ws.Range("I16").Value = 2.98
$ws.Range("L16").Value = 1.23
$ws.Range("Q16").Value = 1.41
$ws.Range("R16").Value = 1.94
$ws.Range("V16").Value = 1.5
$ws.Range("W16").Value = 1.72
$ws.Range("AA16").Value = 48
$ws.Range("AD16").Value = 14.5
# Row 17
$ws.Range("P17").Value = 3.5
$ws.Range("Q17").Value = 1.37
$ws.Range("R17").Value = 2.02
$ws.Range("S17").Value = 1.93
$ws.Range("V17").Value = 1.06
$ws.Range("W17").Value = 4.8
$ws.Range("X17").Value = 44
$ws.Range("Z17").Value = 160
$ws.Range("AB17").Value = 15
$ws.Range("AF17").Value = 10.5
$ws.Range("AL17").Value = 28
$ws.Range("AN17").Value = 3.15
# Row 18
$ws.Range("J18").Value = 6.2
$ws.Range("S18").Value = 2.54
$ws.Range("W18").Value = 4
$ws.Range("X18").Value = 26
$ws.Range("AB18").Value = 9.800000000000001
$ws.Range("AI18").Value = 130
$ws.Range("AO18").Value = 180
# Row 19
$ws.Range("L19").Value = 1.35
$ws.Range("P19").Value = 2.2
$ws.Range("Q19").Value = 1.81
$ws.Range("R19").Value = 1.47
$ws.Range("T19").Value = 1.66
$ws.Range("W19").Value = 1.48
$ws.Range("X19").Value = 16.5
$ws.Range("AE19").Value = 25
$ws.Range("AI19").Value = 34
$ws.Range("AM19").Value = 70
$ws.Range("AN19").Value = 24
# Row 20
$ws.Range("AD20").Value = 11
# Row 21
$ws.Range("F21").Value = 2.38
$ws.Range("G21").Value = 2.42
$ws.Range("H21").Value = 3.4
$ws.Range("I21").Value = 3.5
$ws.Range("L21").Value = 1.48
$ws.Range("N21").Value = 3.25
$ws.Range("P21").Value = 1.76
$ws.Range("S21").Value = 4.2
$ws.Range("U21").Value = 2
$ws.Range("V21").Value = 1.4
$ws.Range("W21").Value = 1.7
$ws.Range("Y21").Value = 11.5
$ws.Range("AB21").Value = 8.800000000000001
$ws.Range("AD21").Value = 15
$ws.Range("AE21").Value = 44
$ws.Range("AF21").Value = 14
$ws.Range("AL21").Value = 48
$ws.Range("AN21").Value = 1000
$ws.Range("AO21").Value = 55
# Row 22
$ws.Range("H22").Value = 1.91
$ws.Range("I22").Value = 1.98
$ws.Range("L22").Value = 1.52
$ws.Range("O22").Value = 1.4
$ws.Range("Q22").Value = 2.16
$ws.Range("Y22").Value = 7.8
$ws.Range("AD22").Value = 10.5
$ws.Range("AE22").Value = 22
# Row 23
$ws.Range("L23").Value = 1.01
$ws.Range("AC23").Value = 7.4
